$d = $word.ActiveDocument

$replacements = @(
    @("2025-12-20 Saturday", "2025-12-21 Sunday"),
    @("651÷6=", "835÷2="),
    @("553÷7=", "912÷9="),
    @("973÷7=", "951÷4="),
    @("206÷6=", "750÷4="),
    @("219÷2=", "899÷3="),
    @("469÷9=", "149÷7="),
    @("905÷2=", "700÷7="),
    @("360÷4=", "112÷6="),
    @("793÷5=", "398÷6="),
    @("532÷8=", "587÷5="),
    @("625÷9=", "242÷5="),
    @("962÷7=", "991÷7="),
    @("751÷8=", "284÷9="),
    @("863÷2=", "804÷5="),
    @("682÷3=", "158÷6="),
    @("474÷7=", "114÷5="),
    @("251÷4=", "472÷3="),
    @("908÷2=", "834÷2="),
    @("647÷9=", "429÷4="),
    @("988÷2=", "732÷7="),
    @("170÷3=", "743÷9="),
    @("924÷4=", "416÷3="),
    @("642÷6=", "427÷4="),
    @("574÷9=", "509÷5="),
    @("527÷7=", "288÷7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
